# Update "想去人数" (want-to-go count) values for several漫展 events.
# These counts are duplicated across sheets because the workbook keeps
# per-category sheets ("展览", "本地生活") as well as a combined
# "全部类型" sheet listing every event.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F8").Value = 239    # 238 -> 239
$wsExhibit.Range("F10").Value = 2310  # 2309 -> 2310
$wsExhibit.Range("F13").Value = 752   # 751 -> 752
$wsExhibit.Range("F18").Value = 24    # 23 -> 24

# 本地生活 sheet
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 2113     # 2112 -> 2113
$wsLocal.Range("F6").Value = 9        # 7 -> 9

# 全部类型 sheet (combined listing)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2113       # 2112 -> 2113
$wsAll.Range("F16").Value = 9         # 7 -> 9
$wsAll.Range("F20").Value = 239       # 238 -> 239
$wsAll.Range("F25").Value = 2310      # 2309 -> 2310
$wsAll.Range("F29").Value = 752       # 751 -> 752
$wsAll.Range("F38").Value = 24        # 23 -> 24
